# Auto-generated: apply scheduled-runner market-data updates to H:N columns
# across 33 leve rows spanning all 8 sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1023.8571
$ws.Range("I129").Value = 650
$ws.Range("J129").Value = 1044.3424
$ws.Range("K129").Value = 1950
$ws.Range("L129").Value = 3133.0272
$ws.Range("M129").Value = 3050
$ws.Range("N129").Value = -13133.0272

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12020.267
$ws.Range("I32").Value = 9165.200000000001
$ws.Range("J32").Value = 34860.8
$ws.Range("K32").Value = 9165.200000000001
$ws.Range("L32").Value = 34860.8
$ws.Range("M32").Value = -8878.200000000001
$ws.Range("N32").Value = -35434.8

$ws.Range("H61").Value = 2274.9524
$ws.Range("I61").Value = 2080.25
$ws.Range("K61").Value = 2080.25
$ws.Range("M61").Value = -1868.25

$ws.Range("H136").Value = 2274.9524
$ws.Range("I136").Value = 2080.25
$ws.Range("K136").Value = 6240.75
$ws.Range("M136").Value = -3690.75

$ws.Range("H139").Value = 67425
$ws.Range("J139").Value = 67425
$ws.Range("L139").Value = 67425
$ws.Range("N139").Value = -77705

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 200
$ws.Range("K8").Value = 200
$ws.Range("M8").Value = -60

$ws.Range("H52").Value = 31988
$ws.Range("J52").Value = 31988
$ws.Range("L52").Value = 31988
$ws.Range("N52").Value = -32514

$ws.Range("H107").Value = 999.8
$ws.Range("I107").Value = 999.75
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 999.75
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 920.25
$ws.Range("N107").Value = -4840

$ws.Range("H121").Value = 31988
$ws.Range("J121").Value = 31988
$ws.Range("L121").Value = 31988
$ws.Range("N121").Value = -35482

$ws.Range("H132").Value = 45837.5
$ws.Range("J132").Value = 45837.5
$ws.Range("L132").Value = 45837.5
$ws.Range("N132").Value = -55957.5

$ws.Range("H134").Value = 1699.3
$ws.Range("I134").Value = 1269.8235
$ws.Range("J134").Value = 2611.9375
$ws.Range("K134").Value = 3809.4705
$ws.Range("L134").Value = 7835.8125
$ws.Range("M134").Value = -1274.4705
$ws.Range("N134").Value = -12905.8125

$ws.Range("H140").Value = 58048.332
$ws.Range("J140").Value = 58048.332
$ws.Range("L140").Value = 58048.332
$ws.Range("N140").Value = -68408.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2710.28
$ws.Range("I31").Value = 1983.3802
$ws.Range("J31").Value = 4489.931
$ws.Range("K31").Value = 1983.3802
$ws.Range("L31").Value = 4489.931
$ws.Range("M31").Value = -1688.3802
$ws.Range("N31").Value = -5079.931

$ws.Range("H34").Value = 2710.28
$ws.Range("I34").Value = 1983.3802
$ws.Range("J34").Value = 4489.931
$ws.Range("K34").Value = 1983.3802
$ws.Range("L34").Value = 4489.931
$ws.Range("M34").Value = -1781.3802
$ws.Range("N34").Value = -4893.931

$ws.Range("H86").Value = 11114279
$ws.Range("I86").Value = 18521640
$ws.Range("J86").Value = 3238.3333
$ws.Range("K86").Value = 18521640
$ws.Range("L86").Value = 3238.3333
$ws.Range("M86").Value = -18520517
$ws.Range("N86").Value = -5484.3333

$ws.Range("H89").Value = 11114279
$ws.Range("I89").Value = 18521640
$ws.Range("J89").Value = 3238.3333
$ws.Range("K89").Value = 92608200
$ws.Range("L89").Value = 16191.6665
$ws.Range("M89").Value = -92602584
$ws.Range("N89").Value = -27423.6665

$ws.Range("H94").Value = 1678.826
$ws.Range("I94").Value = 1976.3846
$ws.Range("J94").Value = 1292
$ws.Range("K94").Value = 1976.3846
$ws.Range("L94").Value = 1292
$ws.Range("M94").Value = -1525.3846
$ws.Range("N94").Value = -2194

$ws.Range("H107").Value = 423405.53
$ws.Range("I107").Value = 43718.74
$ws.Range("J107").Value = 3334337.8
$ws.Range("K107").Value = 43718.74
$ws.Range("L107").Value = 3334337.8
$ws.Range("M107").Value = -41798.74
$ws.Range("N107").Value = -3338177.8

$ws.Range("H140").Value = 88817.664
$ws.Range("J140").Value = 88817.664
$ws.Range("L140").Value = 88817.664
$ws.Range("N140").Value = -99177.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 79825
$ws.Range("J37").Value = 79825
$ws.Range("L37").Value = 239475
$ws.Range("N37").Value = -239699

$ws.Range("H131").Value = 780.1900000000001
$ws.Range("I131").Value = 405.33334
$ws.Range("J131").Value = 846.3412
$ws.Range("K131").Value = 1216.00002
$ws.Range("L131").Value = 2539.0236
$ws.Range("M131").Value = 3823.99998
$ws.Range("N131").Value = -12619.0236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H132").Value = 2453.1765
$ws.Range("I132").Value = 2087.0435
$ws.Range("J132").Value = 3218.7273
$ws.Range("K132").Value = 6261.130500000001
$ws.Range("L132").Value = 9656.1819
$ws.Range("M132").Value = -3731.130500000001
$ws.Range("N132").Value = -14716.1819

$ws.Range("H135").Value = 41963.332
$ws.Range("J135").Value = 41963.332
$ws.Range("L135").Value = 41963.332
$ws.Range("N135").Value = -52103.332

$ws.Range("H138").Value = 50825.3
$ws.Range("J138").Value = 50825.3
$ws.Range("L138").Value = 50825.3
$ws.Range("N138").Value = -61105.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58871.39
$ws.Range("I7").Value = 69645.664
$ws.Range("K7").Value = 69645.664
$ws.Range("M7").Value = -69533.664

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H126").Value = 58871.39
$ws.Range("I126").Value = 69645.664
$ws.Range("K126").Value = 208936.992
$ws.Range("M126").Value = -206466.992

$ws.Range("H127").Value = 57188
$ws.Range("J127").Value = 57188
$ws.Range("L127").Value = 57188
$ws.Range("N127").Value = -67108

$ws.Range("H136").Value = 2709
$ws.Range("I136").Value = 2568.7778
$ws.Range("J136").Value = 3129.6667
$ws.Range("K136").Value = 7706.3334
$ws.Range("L136").Value = 9389.000100000001
$ws.Range("M136").Value = -5156.3334
$ws.Range("N136").Value = -14489.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2177.4727
$ws.Range("I132").Value = 1859.3636
$ws.Range("J132").Value = 2654.6365
$ws.Range("K132").Value = 5578.0908
$ws.Range("L132").Value = 7963.9095
$ws.Range("M132").Value = -3048.0908
$ws.Range("N132").Value = -13023.9095

$ws.Range("H136").Value = 4924.2383
$ws.Range("I136").Value = 3144.4546
$ws.Range("J136").Value = 6882
$ws.Range("K136").Value = 9433.363799999999
$ws.Range("L136").Value = 20646
$ws.Range("M136").Value = -6883.363799999999
$ws.Range("N136").Value = -25746

$ws.Range("H137").Value = 56354.332
$ws.Range("J137").Value = 56354.332
$ws.Range("L137").Value = 56354.332
$ws.Range("N137").Value = -66554.33199999999
